$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to text format so numeric-looking strings are preserved as text,
# matching the original inlineStr cell type used throughout the sheet.
$dRange = $ws.Range("D2:D51")
$dRange.NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = "65.094.50"
$ws.Range("E2").Value = "  +0.24%  "

# Row 3
$ws.Range("D3").Value = "3.183.70"
$ws.Range("E3").Value = "  -1.59%  "

# Row 4
$ws.Range("E4").Value = "  -0.06%  "

# Row 5
$ws.Range("D5").Value = "574.46"
$ws.Range("E5").Value = "  -0.80%  "

# Row 6
$ws.Range("D6").Value = "167.11"
$ws.Range("E6").Value = "  -3.32%  "

# Row 7
$ws.Range("D7").Value = "0.595"
$ws.Range("E7").Value = "  -5.98%  "

# Row 8
$ws.Range("E8").Value = "  -0.03%  "

# Row 9
$ws.Range("D9").Value = "0.119"
$ws.Range("E9").Value = "  -3.10%  "

# Row 10
$ws.Range("D10").Value = "6.69"
$ws.Range("E10").Value = "  -1.43%  "

# Row 11
$ws.Range("D11").Value = "0.387"
$ws.Range("E11").Value = "  -0.84%  "

# Row 12
$ws.Range("D12").Value = "3.737.11"
$ws.Range("E12").Value = "  -1.70%  "

# Row 13
$ws.Range("E13").Value = "  -0.53%  "

# Row 14
$ws.Range("D14").Value = "64.864.22"
$ws.Range("E14").Value = "  -0.32%  "

# Row 15
$ws.Range("D15").Value = "25.57"
$ws.Range("E15").Value = "  -0.79%  "

# Row 16
$ws.Range("D16").Value = "3.183.36"
$ws.Range("E16").Value = "  -1.55%  "

# Row 17
$ws.Range("D17").Value = "0.0000157"
$ws.Range("E17").Value = "  -1.48%  "

# Row 18
$ws.Range("D18").Value = "412.46"
$ws.Range("E18").Value = "  -1.58%  "

# Row 19
$ws.Range("D19").Value = "12.79"
$ws.Range("E19").Value = "  -0.51%  "

# Row 20
$ws.Range("D20").Value = "5.30"
$ws.Range("E20").Value = "  -1.97%  "

# Row 21
$ws.Range("D21").Value = "7.14"
$ws.Range("E21").Value = "  -1.04%  "

# Row 22
$ws.Range("E22").Value = "  +0.14%  "

# Row 23
$ws.Range("D23").Value = "68.95"
$ws.Range("E23").Value = "  -2.60%  "

# Row 24
$ws.Range("E24").Value = "  -1.89%  "

# Row 25
$ws.Range("D25").Value = "0.487"
$ws.Range("E25").Value = "  -1.92%  "

# Row 26
$ws.Range("D26").Value = "0.0000105"
$ws.Range("E26").Value = "  -5.57%  "

# Row 27
$ws.Range("D27").Value = "8.87"
$ws.Range("E27").Value = "  -3.21%  "

# Row 28
$ws.Range("D28").Value = "1.00"
$ws.Range("E28").Value = "  -0.05%  "

# Row 29
$ws.Range("E29").Value = "  -2.70%  "

# Row 30
$ws.Range("D30").Value = "21.43"
$ws.Range("E30").Value = "  -2.33%  "

# Row 31
$ws.Range("D31").Value = "4.98"
$ws.Range("E31").Value = "  -0.92%  "

# Row 32
$ws.Range("D32").Value = "6.35"
$ws.Range("E32").Value = "  -1.49%  "

# Row 33
$ws.Range("E33").Value = "  -2.39%  "

# Row 34
$ws.Range("D34").Value = "155.75"
$ws.Range("E34").Value = "  -0.84%  "

# Row 35
$ws.Range("E35").Value = "  -2.75%  "

# Row 36
$ws.Range("D36").Value = "2.742.52"
$ws.Range("E36").Value = "  -3.33%  "

# Row 37
$ws.Range("D37").Value = "1.73"
$ws.Range("E37").Value = "  -1.33%  "

# Row 38
$ws.Range("D38").Value = "24.01"
$ws.Range("E38").Value = "  -6.09%  "

# Row 39
$ws.Range("D39").Value = "4.13"
$ws.Range("E39").Value = "  -2.90%  "

# Row 40
$ws.Range("E40").Value = "  -2.37%  "

# Row 41
$ws.Range("D41").Value = "0.0634"
$ws.Range("E41").Value = "  +0.31%  "

# Row 42
$ws.Range("D42").Value = "5.54"
$ws.Range("E42").Value = "  -4.04%  "

# Row 43
$ws.Range("E43").Value = "  -1.38%  "

# Row 44
$ws.Range("D44").Value = "293.54"
$ws.Range("E44").Value = "  -3.53%  "

# Row 45
$ws.Range("D45").Value = "21.41"
$ws.Range("E45").Value = "  -3.82%  "

# Row 46
$ws.Range("D46").Value = "0.999"
$ws.Range("E46").Value = "  -0.08%  "

# Row 47
$ws.Range("D47").Value = "0.0988"
$ws.Range("E47").Value = "  -2.79%  "

# Row 48
$ws.Range("D48").Value = "1.97"
$ws.Range("E48").Value = "  -9.71%  "

# Row 49
$ws.Range("B49").Value = "Cosmos"
$ws.Range("C49").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D49").Value = "5.79"
$ws.Range("E49").Value = "  -1.04%  "

# Row 50
$ws.Range("B50").Value = "WhiteBITCoin"
$ws.Range("C50").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D50").Value = "10.44"
$ws.Range("E50").Value = "  +0.41%  "

# Row 51
$ws.Range("D51").Value = "0.900"
$ws.Range("E51").Value = "  -3.83%  "

# Remove the temporary number formatting so styling matches the original (no style index).
$dRange.ClearFormats()
